# Refresh the cryptos table with the latest coinranking.com scrape:
# updated prices / 1h volume-change percentages, plus a handful of
# coins that swapped rank (and therefore row) with their neighbour.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.550.15'
$ws.Range("E2").Value = '''  +6.65%  '
$ws.Range("D3").Value = '''2.029.64'
$ws.Range("E3").Value = '''  +7.80%  '
$ws.Range("E4").Value = '''  +0.08%  '
$ws.Range("D5").Value = '''254.04'
$ws.Range("E5").Value = '''  +3.84%  '
$ws.Range("D6").Value = '''0.696'
$ws.Range("E6").Value = '''  +1.47%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '''  +0.04%  '
$ws.Range("D8").Value = '''46.71'
$ws.Range("E8").Value = '''  +9.88%  '
$ws.Range("D9").Value = '''0.383'
$ws.Range("E9").Value = '''  +8.75%  '
$ws.Range("D10").Value = '''58.10'
$ws.Range("E10").Value = '''  +5.56%  '
$ws.Range("D11").Value = '''0.0772'
$ws.Range("E11").Value = '''  +4.65%  '
$ws.Range("E12").Value = '''  +2.61%  '
$ws.Range("D13").Value = '''15.54'
$ws.Range("E13").Value = '''  +13.43%  '
$ws.Range("D14").Value = '''0.843'
$ws.Range("E14").Value = '''  +8.38%  '
$ws.Range("D15").Value = '''2.338.01'
$ws.Range("E15").Value = '''  +8.33%  '
$ws.Range("D16").Value = '''5.22'
$ws.Range("E16").Value = '''  +5.62%  '
$ws.Range("D17").Value = '''2.036.12'
$ws.Range("E17").Value = '''  +8.21%  '
$ws.Range("D18").Value = '''37.589.95'
$ws.Range("E18").Value = '''  +6.60%  '
$ws.Range("D19").Value = '''75.67'
$ws.Range("E19").Value = '''  +3.63%  '
$ws.Range("D20").Value = '''0.0₃0864'
$ws.Range("E20").Value = '''  +5.69%  '
$ws.Range("D21").Value = '''13.79'
$ws.Range("E21").Value = '''  +8.39%  '
$ws.Range("D22").Value = '''255.68'
$ws.Range("E22").Value = '''  +5.13%  '
$ws.Range("D23").Value = '''5.30'
$ws.Range("E23").Value = '''  +3.64%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '''  -0.05%  '
$ws.Range("D25").Value = '''2.53'
$ws.Range("E25").Value = '''  -4.04%  '
$ws.Range("D26").Value = '''169.97'
$ws.Range("E26").Value = '''  +2.09%  '
$ws.Range("D27").Value = '''2.16'
$ws.Range("E27").Value = '''  +1.23%  '
$ws.Range("D28").Value = '''8.94'
$ws.Range("E28").Value = '''  +5.67%  '
$ws.Range("D29").Value = '''20.22'
$ws.Range("E29").Value = '''  +11.13%  '
$ws.Range("E30").Value = '''  +2.60%  '
$ws.Range("D31").Value = '''4.61'
$ws.Range("E31").Value = '''  +8.18%  '
$ws.Range("B32").Value = '''Gas'
$ws.Range("C32").Value = '''https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D32").Value = '''22.16'
$ws.Range("E32").Value = '''  +68.17%  '
$ws.Range("B33").Value = '''Hedera'
$ws.Range("C33").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.0617'
$ws.Range("E33").Value = '''  +4.47%  '
$ws.Range("D34").Value = '''0.0917'
$ws.Range("E34").Value = '''  +28.12%  '
$ws.Range("D35").Value = '''4.37'
$ws.Range("E35").Value = '''  +5.58%  '
$ws.Range("B36").Value = '''BinanceUSD'
$ws.Range("C36").Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '''  +0.08%  '
$ws.Range("B37").Value = '''WEMIXToken'
$ws.Range("C37").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '''1.88'
$ws.Range("E37").Value = '''  +0.17%  '
$ws.Range("D38").Value = '''0.897'
$ws.Range("E38").Value = '''  +6.61%  '
$ws.Range("B39").Value = '''LidoDAOToken'
$ws.Range("C39").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.20'
$ws.Range("E39").Value = '''  +14.58%  '
$ws.Range("B40").Value = '''TrustWalletToken'
$ws.Range("C40").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.46'
$ws.Range("E40").Value = '''  +0.65%  '
$ws.Range("D41").Value = '''103.74'
$ws.Range("E41").Value = '''  +6.58%  '
$ws.Range("E42").Value = '''  +4.59%  '
$ws.Range("D43").Value = '''17.50'
$ws.Range("E43").Value = '''  +2.77%  '
$ws.Range("D44").Value = '''1.13'
$ws.Range("E44").Value = '''  +6.25%  '
$ws.Range("D45").Value = '''2.89'
$ws.Range("E45").Value = '''  +20.50%  '
$ws.Range("D46").Value = '''1.373.41'
$ws.Range("E46").Value = '''  +3.85%  '
$ws.Range("B47").Value = '''Cronos'
$ws.Range("C47").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.0851'
$ws.Range("E47").Value = '''  +5.44%  '
$ws.Range("B48").Value = '''RenderToken'
$ws.Range("C48").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''2.42'
$ws.Range("E48").Value = '''  +3.29%  '
$ws.Range("E49").Value = '''  +5.04%  '
$ws.Range("D50").Value = '''3.94'
$ws.Range("E50").Value = '''  +19.51%  '
$ws.Range("D51").Value = '''6.61'
$ws.Range("E51").Value = '''  +5.80%  '
